$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "section" column (D) for the
# new "course" field, shifting section/batch/... one column to the right.
$ws.Columns("D").Insert()

# New column D header + sample/help value.
$ws.Range("D1").Value = "course"
$ws.Range("D2").Value = "1,2,3 (depending on MBA,BBA,Btech respectively)"

# Match the column width used for the new column.
$ws.Columns("D").ColumnWidth = 45.166666666666664

# The student_email_id hyperlink lived on the old I2; after the column
# insert its data moved to J2, so re-anchor the hyperlink there.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:777@dd.com") | Out-Null

# Restore the selection Excel left the workbook in after the edit.
$ws.Range("D6").Select() | Out-Null
